$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 13688.556
$ws.Range("I11").Value = 13688.556
$ws.Range("K11").Value = 13688.556
$ws.Range("M11").Value = -13548.556
$ws.Range("H19").Value = 493.33334
$ws.Range("J19").Value = 595
$ws.Range("L19").Value = 595
$ws.Range("N19").Value = -945
$ws.Range("H41").Value = 682.8
$ws.Range("I41").Value = 958
$ws.Range("J41").Value = 270
$ws.Range("K41").Value = 958
$ws.Range("L41").Value = 270
$ws.Range("M41").Value = -518
$ws.Range("N41").Value = -1150
$ws.Range("H113").Value = 5999.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5999.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5999.5
$ws.Range("N113").Value = -12507.5
$ws.Range("H132").Value = 35231.332
$ws.Range("I132").Value = 40466.69
$ws.Range("K132").Value = 121400.07
$ws.Range("M132").Value = -118870.07
$ws.Range("H135").Value = 1211.7142
$ws.Range("I135").Value = 1211.7142
$ws.Range("K135").Value = 10905.4278
$ws.Range("M135").Value = -8370.427799999999
$ws.Range("H137").Value = 9094417
$ws.Range("I137").Value = 3148.125
$ws.Range("K137").Value = 9444.375
$ws.Range("M137").Value = -6894.375
$ws.Range("H138").Value = 5602.2812
$ws.Range("I138").Value = 6750.6553
$ws.Range("J138").Value = 4650.7715
$ws.Range("K138").Value = 20251.9659
$ws.Range("L138").Value = 13952.3145
$ws.Range("M138").Value = -15111.9659
$ws.Range("N138").Value = -24232.3145
$ws.Range("H141").Value = 1752.2778
$ws.Range("I141").Value = 1213.2069
$ws.Range("J141").Value = 3985.5715
$ws.Range("K141").Value = 3639.620699999999
$ws.Range("L141").Value = 11956.7145
$ws.Range("M141").Value = 1540.379300000001
$ws.Range("N141").Value = -22316.7145
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 107356.625
$ws.Range("I32").Value = 110533.13
$ws.Range("K32").Value = 110533.13
$ws.Range("M32").Value = -110246.13
$ws.Range("H45").Value = 73834.78999999999
$ws.Range("I45").Value = 93162.45
$ws.Range("K45").Value = 93162.45
$ws.Range("M45").Value = -92785.45
$ws.Range("H46").Value = 11481.125
$ws.Range("I46").Value = 3999
$ws.Range("J46").Value = 12550
$ws.Range("K46").Value = 3999
$ws.Range("L46").Value = 12550
$ws.Range("M46").Value = -3680
$ws.Range("N46").Value = -13188
$ws.Range("H74").Value = 638900.75
$ws.Range("I74").Value = 1379
$ws.Range("J74").Value = 1241004.6
$ws.Range("K74").Value = 1379
$ws.Range("L74").Value = 1241004.6
$ws.Range("M74").Value = -505
$ws.Range("N74").Value = -1242752.6
$ws.Range("H77").Value = 638900.75
$ws.Range("I77").Value = 1379
$ws.Range("J77").Value = 1241004.6
$ws.Range("K77").Value = 6895
$ws.Range("L77").Value = 6205023
$ws.Range("M77").Value = -2527
$ws.Range("N77").Value = -6213759
$ws.Range("H97").Value = 8071.933
$ws.Range("I97").Value = 9848.25
$ws.Range("J97").Value = 966.6667
$ws.Range("K97").Value = 9848.25
$ws.Range("L97").Value = 966.6667
$ws.Range("M97").Value = -9352.25
$ws.Range("N97").Value = -1958.6667
$ws.Range("H102").Value = 1661.5769
$ws.Range("I102").Value = 1618.04
$ws.Range("J102").Value = 2750
$ws.Range("K102").Value = 1618.04
$ws.Range("L102").Value = 2750
$ws.Range("M102").Value = 3.960000000000036
$ws.Range("N102").Value = -5994
$ws.Range("H122").Value = 1599.4615
$ws.Range("I122").Value = 1686.9131
$ws.Range("J122").Value = 929
$ws.Range("K122").Value = 5060.7393
$ws.Range("L122").Value = 2787
$ws.Range("M122").Value = -2610.7393
$ws.Range("N122").Value = -7687

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1141.3103
$ws.Range("I20").Value = 943.7273
$ws.Range("J20").Value = 1262.0555
$ws.Range("K20").Value = 943.7273
$ws.Range("L20").Value = 1262.0555
$ws.Range("M20").Value = -696.7273
$ws.Range("N20").Value = -1756.0555
$ws.Range("H81").Value = 44593.332
$ws.Range("J81").Value = 44593.332
$ws.Range("L81").Value = 44593.332
$ws.Range("N81").Value = -46715.332
$ws.Range("H84").Value = 44593.332
$ws.Range("J84").Value = 44593.332
$ws.Range("L84").Value = 133779.996
$ws.Range("N84").Value = -144387.996
$ws.Range("H86").Value = 1967.96
$ws.Range("I86").Value = 1176.7646
$ws.Range("J86").Value = 3649.25
$ws.Range("K86").Value = 1176.7646
$ws.Range("L86").Value = 3649.25
$ws.Range("M86").Value = -53.76459999999997
$ws.Range("N86").Value = -5895.25
$ws.Range("H89").Value = 1967.96
$ws.Range("I89").Value = 1176.7646
$ws.Range("J89").Value = 3649.25
$ws.Range("K89").Value = 5883.823
$ws.Range("L89").Value = 18246.25
$ws.Range("M89").Value = -267.8230000000003
$ws.Range("N89").Value = -29478.25
$ws.Range("H99").Value = 7270.625
$ws.Range("I99").Value = 11215.077
$ws.Range("K99").Value = 11215.077
$ws.Range("M99").Value = -9717.076999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 35000
$ws.Range("J43").Value = 35000
$ws.Range("L43").Value = 35000
$ws.Range("N43").Value = -35368
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("H101").Value = 35000
$ws.Range("J101").Value = 35000
$ws.Range("L101").Value = 35000
$ws.Range("N101").Value = -41490
$ws.Range("H107").Value = 1410.5
$ws.Range("I107").Value = 1284.8
$ws.Range("J107").Value = 1620
$ws.Range("K107").Value = 1284.8
$ws.Range("L107").Value = 1620
$ws.Range("M107").Value = 635.2
$ws.Range("N107").Value = -5460
$ws.Range("H122").Value = 2806.2942
$ws.Range("I122").Value = 2948.9
$ws.Range("J122").Value = 2602.5715
$ws.Range("K122").Value = 8846.700000000001
$ws.Range("L122").Value = 7807.7145
$ws.Range("M122").Value = -6396.700000000001
$ws.Range("N122").Value = -12707.7145
$ws.Range("H132").Value = 2520.3396
$ws.Range("I132").Value = 2348.1191
$ws.Range("J132").Value = 3177.9092
$ws.Range("K132").Value = 7044.3573
$ws.Range("L132").Value = 9533.7276
$ws.Range("M132").Value = -4514.3573
$ws.Range("N132").Value = -14593.7276
$ws.Range("H134").Value = 1724.766
$ws.Range("I134").Value = 1186.4412
$ws.Range("J134").Value = 3132.6924
$ws.Range("K134").Value = 3559.3236
$ws.Range("L134").Value = 9398.0772
$ws.Range("M134").Value = -1024.3236
$ws.Range("N134").Value = -14468.0772
$ws.Range("N88").ClearContents()
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1958.7
$ws.Range("I86").Value = 944.5714
$ws.Range("J86").Value = 4325
$ws.Range("K86").Value = 2833.7142
$ws.Range("L86").Value = 12975
$ws.Range("M86").Value = -1647.7142
$ws.Range("N86").Value = -15347
$ws.Range("H89").Value = 1958.7
$ws.Range("I89").Value = 944.5714
$ws.Range("J89").Value = 4325
$ws.Range("K89").Value = 8501.142600000001
$ws.Range("L89").Value = 38925
$ws.Range("M89").Value = -2573.142600000001
$ws.Range("N89").Value = -50781
$ws.Range("H131").Value = 3406277.8
$ws.Range("I131").Value = 8266118.5
$ws.Range("J131").Value = 65137.312
$ws.Range("K131").Value = 24798355.5
$ws.Range("L131").Value = 195411.936
$ws.Range("M131").Value = -24793315.5
$ws.Range("N131").Value = -205491.936

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 4000000
$ws.Range("I20").Value = 4000000
$ws.Range("K20").Value = 4000000
$ws.Range("M20").Value = -3999755
$ws.Range("H132").Value = 560040.2
$ws.Range("I132").Value = 2207.1892
$ws.Range("K132").Value = 6621.567599999999
$ws.Range("M132").Value = -4091.567599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1294384.1
$ws.Range("I40").Value = 1545133.4
$ws.Range("J40").Value = 4816.2856
$ws.Range("K40").Value = 1545133.4
$ws.Range("L40").Value = 4816.2856
$ws.Range("M40").Value = -1544997.4
$ws.Range("N40").Value = -5088.2856
$ws.Range("H93").Value = 1975.72
$ws.Range("I93").Value = 2018.762
$ws.Range("K93").Value = 2018.762
$ws.Range("M93").Value = -770.7619999999999
$ws.Range("H122").Value = 3899.6155
$ws.Range("I122").Value = 3280.5
$ws.Range("K122").Value = 9841.5
$ws.Range("M122").Value = -7391.5
$ws.Range("H132").Value = 7873.773
$ws.Range("I132").Value = 2301.6428
$ws.Range("J132").Value = 17625
$ws.Range("K132").Value = 6904.928400000001
$ws.Range("L132").Value = 52875
$ws.Range("M132").Value = -4374.928400000001
$ws.Range("N132").Value = -57935
$ws.Range("H136").Value = 4501
$ws.Range("I136").Value = 1023
$ws.Range("K136").Value = 3069
$ws.Range("M136").Value = -519

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 676689.75
$ws.Range("I100").Value = 934202.9
$ws.Range("J100").Value = 717.75
$ws.Range("K100").Value = 1868405.8
$ws.Range("L100").Value = 1435.5
$ws.Range("M100").Value = -1867864.8
$ws.Range("N100").Value = -2517.5
$ws.Range("H132").Value = 23727.979
$ws.Range("I132").Value = 29416
$ws.Range("K132").Value = 88248
$ws.Range("M132").Value = -85718
$ws.Range("H136").Value = 17687.08
$ws.Range("I136").Value = 24888.643
$ws.Range("J136").Value = 3283.9524
$ws.Range("K136").Value = 74665.929
$ws.Range("L136").Value = 9851.8572
$ws.Range("M136").Value = -72115.929
$ws.Range("N136").Value = -14951.8572

Write-Host "Applied all cell updates"